$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Test_01"
$ws.Range("A3").Value = "Test_02"
$ws.Range("A4").Select()
